# Refresh the cryptos list: price (D) and 1h change (E) for every coin row,
# and swap the OKB / Filecoin rows (41 <-> 42) to match the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column cells hold plain text in this sheet (e.g. '60.496.87' uses dots
# as thousands separators). Force text format before assigning so Excel's COM
# layer doesn't auto-coerce the digit-looking strings into numbers, then drop
# back to the sheet's normal (unstyled) look so formatting matches the rest of
# the table.
function Set-TextValue($range, $value) {
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

# Row 2
Set-TextValue $ws.Range("D2") '60.496.87'
$ws.Range("E2").Value = '  -0.49%  '
# Row 3
Set-TextValue $ws.Range("D3") '2.579.26'
$ws.Range("E3").Value = '  -3.95%  '
# Row 4
$ws.Range("E4").Value = '  -0.31%  '
# Row 5
Set-TextValue $ws.Range("D5") '507.91'
$ws.Range("E5").Value = '  -1.46%  '
# Row 6
Set-TextValue $ws.Range("D6") '155.77'
$ws.Range("E6").Value = '  -4.00%  '
# Row 7
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  +0.02%  '
# Row 8
$ws.Range("E8").Value = '  -5.27%  '
# Row 9
Set-TextValue $ws.Range("D9") '2.588.42'
$ws.Range("E9").Value = '  -3.77%  '
# Row 10
Set-TextValue $ws.Range("D10") '6.67'
$ws.Range("E10").Value = '  +8.22%  '
# Row 11
$ws.Range("E11").Value = '  -3.24%  '
# Row 12
Set-TextValue $ws.Range("D12") '0.346'
$ws.Range("E12").Value = '  -1.51%  '
# Row 13
Set-TextValue $ws.Range("D13") '0.129'
$ws.Range("E13").Value = '  +1.23%  '
# Row 14
Set-TextValue $ws.Range("D14") '3.041.67'
$ws.Range("E14").Value = '  -3.07%  '
# Row 15
Set-TextValue $ws.Range("D15") '60.532.57'
$ws.Range("E15").Value = '  -1.05%  '
# Row 16
Set-TextValue $ws.Range("D16") '21.68'
$ws.Range("E16").Value = '  -4.45%  '
# Row 17
$ws.Range("E17").Value = '  -1.27%  '
# Row 18
Set-TextValue $ws.Range("D18") '2.592.31'
$ws.Range("E18").Value = '  -3.54%  '
# Row 19
Set-TextValue $ws.Range("D19") '4.77'
$ws.Range("E19").Value = '  -2.54%  '
# Row 20
Set-TextValue $ws.Range("D20") '346.76'
$ws.Range("E20").Value = '  -2.52%  '
# Row 21
Set-TextValue $ws.Range("D21") '10.50'
$ws.Range("E21").Value = '  -1.46%  '
# Row 22
Set-TextValue $ws.Range("D22") '6.12'
$ws.Range("E22").Value = '  -1.96%  '
# Row 23
$ws.Range("E23").Value = '  -0.06%  '
# Row 24
Set-TextValue $ws.Range("D24") '60.65'
$ws.Range("E24").Value = '  -0.09%  '
# Row 25
$ws.Range("E25").Value = '  -2.00%  '
# Row 26
$ws.Range("E26").Value = '  -1.68%  '
# Row 27
Set-TextValue $ws.Range("D27") '2.709.08'
$ws.Range("E27").Value = '  -2.51%  '
# Row 28
Set-TextValue $ws.Range("D28") '0.972'
$ws.Range("E28").Value = '  -2.63%  '
# Row 29
Set-TextValue $ws.Range("D29") '0.0₃0846'
$ws.Range("E29").Value = '  -3.81%  '
# Row 30
Set-TextValue $ws.Range("D30") '7.42'
$ws.Range("E30").Value = '  -3.10%  '
# Row 31
$ws.Range("E31").Value = '  +0.03%  '
# Row 32
Set-TextValue $ws.Range("D32") '19.40'
$ws.Range("E32").Value = '  -2.19%  '
# Row 33
Set-TextValue $ws.Range("D33") '152.95'
$ws.Range("E33").Value = '  -3.76%  '
# Row 34
Set-TextValue $ws.Range("D34") '1.56'
$ws.Range("E34").Value = '  -2.09%  '
# Row 35
Set-TextValue $ws.Range("D35") '5.73'
$ws.Range("E35").Value = '  -0.76%  '
# Row 36
Set-TextValue $ws.Range("D36") '4.02'
$ws.Range("E36").Value = '  -1.90%  '
# Row 37
$ws.Range("E37").Value = '  -3.52%  '
# Row 38
Set-TextValue $ws.Range("D38") '0.852'
$ws.Range("E38").Value = '  +0.08%  '
# Row 39
Set-TextValue $ws.Range("D39") '1.48'
$ws.Range("E39").Value = '  -4.77%  '
# Row 40
Set-TextValue $ws.Range("D40") '0.845'
$ws.Range("E40").Value = '  -4.70%  '
# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D41") '3.77'
$ws.Range("E41").Value = '  -2.18%  '
# Row 42
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D42") '36.11'
$ws.Range("E42").Value = '  -0.05%  '
# Row 43
Set-TextValue $ws.Range("D43") '296.92'
$ws.Range("E43").Value = '  -1.99%  '
# Row 44
$ws.Range("E44").Value = '  -4.29%  '
# Row 45
Set-TextValue $ws.Range("D45") '0.1000'
$ws.Range("E45").Value = '  -1.74%  '
# Row 46
$ws.Range("E46").Value = '  -3.78%  '
# Row 47
Set-TextValue $ws.Range("D47") '0.999'
$ws.Range("E47").Value = '  +0.35%  '
# Row 48
Set-TextValue $ws.Range("D48") '19.77'
$ws.Range("E48").Value = '  -2.49%  '
# Row 49
Set-TextValue $ws.Range("D49") '4.83'
$ws.Range("E49").Value = '  -4.66%  '
# Row 50
Set-TextValue $ws.Range("D50") '0.0234'
$ws.Range("E50").Value = '  -2.90%  '
# Row 51
Set-TextValue $ws.Range("D51") '10.30'
$ws.Range("E51").Value = '  +0.26%  '
